# Update the "Förändrad" (Changed) date column (C) for data rows 2-29
# from 2024-10-13 (serial 45578) to 2024-10-14 (serial 45579).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45578) {
        $cell.Value = 45579
    }
}
